# NIT-9000820550 "Estado de Cuenta" refresh.
# The prior period rows are replaced: data is now grouped by worker
# (YURIS ZAPATEIRO GUZMAN first, then JESUS MARIA DE LA ROSA PEREZ),
# each with periods listed 2009 -> 1904 (descending), rewriting
# C:F for rows 16-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$workers = @(
    [pscustomobject]@{
        Doc  = "1047396267"
        Name = "YURIS ZAPATEIRO GUZMAN"
        Periods = @(
            @{ Period = "2009"; Mora = 29166 },
            @{ Period = "2008"; Mora = 31249 },
            @{ Period = "2007"; Mora = 31249 },
            @{ Period = "2006"; Mora = 31249 },
            @{ Period = "2005"; Mora = 31249 },
            @{ Period = "2004"; Mora = 31249 },
            @{ Period = "2003"; Mora = 31249 },
            @{ Period = "2002"; Mora = 31249 },
            @{ Period = "2001"; Mora = 31249 },
            @{ Period = "1912"; Mora = 31249 },
            @{ Period = "1911"; Mora = 31249 },
            @{ Period = "1910"; Mora = 31249 },
            @{ Period = "1909"; Mora = 31249 },
            @{ Period = "1908"; Mora = 31249 },
            @{ Period = "1907"; Mora = 31249 },
            @{ Period = "1906"; Mora = 31249 },
            @{ Period = "1905"; Mora = 33125 },
            @{ Period = "1904"; Mora = 33125 }
        )
    },
    [pscustomobject]@{
        Doc  = "73213618"
        Name = "JESUS MARIA DE LA ROSA PEREZ"
        Periods = @(
            @{ Period = "2009"; Mora = 30916 },
            @{ Period = "2008"; Mora = 33125 },
            @{ Period = "2007"; Mora = 33125 },
            @{ Period = "2006"; Mora = 33125 },
            @{ Period = "2005"; Mora = 33125 },
            @{ Period = "2004"; Mora = 33125 },
            @{ Period = "2003"; Mora = 33125 },
            @{ Period = "2002"; Mora = 33125 },
            @{ Period = "2001"; Mora = 33125 },
            @{ Period = "1912"; Mora = 33125 },
            @{ Period = "1911"; Mora = 33125 },
            @{ Period = "1910"; Mora = 33125 },
            @{ Period = "1909"; Mora = 33125 },
            @{ Period = "1908"; Mora = 33125 },
            @{ Period = "1907"; Mora = 33125 },
            @{ Period = "1906"; Mora = 33125 },
            @{ Period = "1905"; Mora = 33125 },
            @{ Period = "1904"; Mora = 33125 }
        )
    }
)

$row = 16
foreach ($worker in $workers) {
    foreach ($p in $worker.Periods) {
        $ws.Range("C$row").Value = $worker.Doc
        $ws.Range("D$row").Value = $worker.Name
        $ws.Range("E$row").Value = $p.Period
        $ws.Range("F$row").Value = $p.Mora
        $row = $row + 1
    }
}
